# "more plate from 0 to 1"
# The 0-indexed scanner/lab "Sample" IDs in column D of Sheet1 (e.g. "s3181")
# are normalised to upper-case "Lab ID"s ("S3181"), and the column header
# itself is renamed from "Sample" to "Lab ID".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$lastRow = $ws.Cells.Item($ws.Rows.Count, 4).End(-4162).Row  # xlUp

# Header
$ws.Cells.Item(1, 4).Value = "Lab ID"

# Body: upper-case every existing Sample/Lab ID value in column D
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    $val = $cell.Value2
    if ($null -ne $val) {
        $cell.Value = ([string]$val).ToUpper()
    }
}

# Reflect the view state captured in the diff: the user had scrolled the
# window and selected the whole of column D before saving.
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 88
$ws.Range("D1:D1048576").Select() | Out-Null
